$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was id=11, now id=15)
$ws.Range("A2").Value = 15
$ws.Range("E2").Value = 0.8882575757575757
$ws.Range("F2").Value = 0.9583333333333333
$ws.Range("G2").Value = 0.8181818181818181
$ws.Range("H2").Value = "Q08IB7, P42236, Q1JUP4, Q08IC0"
$ws.Range("I2").Value = "Q3ZCJ2, Q9JII6, Q5R5D5, P50578, P14550, P51635, O70473"

# Row 3 (was id=12, now id=5)
$ws.Range("A3").Value = 5
$ws.Range("E3").Value = 0.7626262626262625
$ws.Range("F3").Value = 0.7626262626262625
$ws.Range("H3").Value = "Q8PLG9, O13471, Q8GKZ1, P62459, Q8X8T3, Q5LL27, Q47AL5, Q5FS83, O30027, P63952, P60859, Q5H0L1, Q492K3, Q5LAZ8, Q7MLS6, Q9L6I1, Q8PZR8, Q9F854, Q8XV79, Q3IQR5, P57201, Q7U8K7, Q7V004, Q318Q4, Q323J2, Q3A133, P62460, Q3MEV7, Q7M8K9, Q5F7D8, P60862, Q7NL02, Q722Y3, Q8DTQ7, Q67KH6, Q9C5U8, Q64RE7, P59400, Q8ABA9, Q5YYQ0, Q3J7H1, Q2YAU5, Q5P790, Q9RH05, Q87C29, Q47QS9, Q5WX91, Q31GY9, P00815, Q3SI67, O59626, Q8CQ95, Q7VQX0, Q5KVC6, P10370, Q8ZY17, Q9PBC5, O33775, Q989E7, Q87WV5, Q5PDP3, P73058, Q845V3, Q6G5Z7, Q465R4, Q3JMZ6, Q4FNE2, Q2P3K3, Q606Q2, Q9HVW9, Q32EE9, Q4JW59, Q8P9P2, Q4A044, P16245, Q3Z878, Q3ZXL7, Q3J079, Q4QN74, Q3BUF7, Q5ZW87, P63951, Q57F99, Q7VA26, Q4FQF8, Q92S26, Q3SW26, P18786, P62458, Q7VFF5, Q8TXG3, Q63Q86, Q62GD9, Q9RQ82, Q8ESR8, Q92E84, Q7N6I1, Q58851, Q8F393, Q2YZB4, Q65RB1, Q5SKC1, Q9HPW5, Q8G2R2, P28736, Q4UU40, Q2RGV8, Q8YSM8, Q7W2Y4, Q8YF59, Q8TL41, Q8R882, Q5HCL9, Q87QL1, Q6FEC8, Q81T62, Q8Y9G1, Q66C49, Q02136, P59398, Q3AD53, Q9CLM4, P62457, Q8KEY6, Q9K6Z2, Q5MZ26, B0SS41, Q9RQ88, Q4ZNV9, Q5HSJ3, Q6ANR3, Q7VSZ1, Q8DGR2, Q97KI2, Q970Y9, Q39K91, Q48EC9, Q47XB8, Q3ICF0, Q81G06, Q2RQM7, Q5E638, Q5WDH8, Q8NNT5, P60858, B0S9F2, P63953, O26327, Q65EF9, Q5HKN8, P60860, Q6GDC6, Q9ZHE6, Q9PM77, Q8EFB1, Q46J19, Q82AA6, Q4J8I8, Q3M5D3, P45353, Q6KZD1, Q9JYH8, Q30TA4, Q57MS3, Q5QWP8, O66976, Q7P0F6, Q39YP7, Q4KI73, P59397, Q5NAY4, O34651, Q6D411, Q9A5V1, Q9RSI4, Q313T7, Q8G4S9, Q9CC57, Q8YWL4, Q8Z5K0, O74712, Q46WL2, Q6HLE7, P62456, Q5V574, P59401, Q9JTH9, Q63DX2, Q8UHX1, Q8ZFX5, Q5JFR5, Q9RQ85, Q5X5W9, Q8D8Q0, P63954, Q9X0D1, Q2YPB8, P24226, P60861, P9WNW8, P44001, P07685, Q7WDY4, P59399, Q7UX39, Q7V5N9, Q8FNZ0, Q82W26, Q8GDP4, Q12670, Q3AR25, Q6AE76, Q6A8L5, P06988, P9WNW9, Q3Z0G5, Q3B5E3, Q3KHZ2, Q8FG52, Q9P777"

# Row 4 (was id=3, now id=14)
$ws.Range("A4").Value = 14
$ws.Range("D4").Value = 57.98770866620097
$ws.Range("E4").Value = 0.8426616915422886
$ws.Range("F4").Value = 0.75
$ws.Range("G4").Value = 0.9353233830845772
$ws.Range("H4").Value = "G2SYC0, B0MC58"
$ws.Range("I4").Value = "Q8WVX9, Q5R834, Q7ZXF5, Q5ZM72, Q922J9, Q66H50"

# Row 5 (was id=15, now id=1)
$ws.Range("A5").Value = 1
$ws.Range("D5").Value = 51.75746525426651
$ws.Range("E5").Value = 0.7026153952459924
$ws.Range("F5").Value = 0.4699074074074074
$ws.Range("G5").Value = 0.9353233830845772
$ws.Range("H5").Value = "Q91VA0, P33124, O06417, Q9QUJ7, P9WQ51, Q63151, O88561, Q2KHW5, Q5R668, Q4R4P9, P18163, Q99PU5, Q7ZYC4, P47912, O88813, P9WQ37, O14975, P97524, O05307, Q8JZR0, O35488, Q91VE0, Q9M0X9, O95573, Q8LPS1, Q5ZKR7, Q9BEA2, Q8LKS5, Q84P21, Q9CZW4, P33121, Q08AH1, Q84P23, P41216, Q5FVE4, Q2XU92, P39518, Q4R3Y4, Q5RDY4, Q91WC3, O53521, O60488, P39002, Q5K4L6, P69451, Q5SKN9, Q9UKU0, Q9ULC5, A1L1K7, O35547, O53551, Q6P1M0, Q924N5, Q96GR2, P30624, Q9JID6"
$ws.Range("I5").Value = "Q8WVX9, Q5R834, Q7ZXF5, Q5ZM72, Q922J9, Q66H50"

# Row 6 (was id=2, now id=4)
$ws.Range("A6").Value = 4
$ws.Range("D6").Value = 15.31924337854084
$ws.Range("E6").Value = 0.647870024875622
$ws.Range("F6").Value = 0.3604166666666667
$ws.Range("G6").Value = 0.9353233830845772
$ws.Range("H6").Value = "P0AGG2, P0AGG3"
$ws.Range("I6").Value = "Q8WVX9, Q5R834, Q7ZXF5, Q5ZM72, Q922J9, Q66H50"
